# Applies the workbook update described in the commit:
# "WIP saving loans , downpayments and cash flows"
#
# Touches six sheets:
#   times                - roll StartTime/StopTime forward by ~4 years
#   scenario_data_emlab  - bump scenario year + Co2/fuel prices
#   conventionals        - replace OIL & LIGNITE plants with two new
#                           NATURAL_GAS placeholder plants, drop two rows
#   renewables           - prepend three new placeholder plants, drop the
#                           RunOfRiver plant
#   storages             - add a new STORAGE placeholder plant (was empty)
#   biogas               - prepend a new Biogas placeholder plant

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# times: StartTime / StopTime shift forward
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("times")
$ws.Cells.Item(2, 2).Value = 45291.99861111111   # StartTime
$ws.Cells.Item(3, 2).Value = 45656.99861111111   # StopTime

# ---------------------------------------------------------------------
# scenario_data_emlab: year + price updates
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("scenario_data_emlab")
$ws.Cells.Item(1, 2).Value = 2024                   # scenario year
$ws.Cells.Item(2, 2).Value = 36.31999999999999      # Co2Prices
$ws.Cells.Item(5, 2).Value = 13.616                 # FuelPrice_HARD_COAL
$ws.Cells.Item(6, 2).Value = 21.392                 # FuelPrice_NATURAL_GAS
$ws.Cells.Item(7, 2).Value = 53.136                 # FuelPrice_OIL

# ---------------------------------------------------------------------
# conventionals: drop OIL + LIGNITE rows, shift NUCLEAR/NATURAL_GAS up,
# and turn the first two rows into new placeholder NATURAL_GAS plants.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("conventionals")
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# row 2: new placeholder NATURAL_GAS plant
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 99991700006
$ws.Cells.Item(2, 3).Value = "NATURAL_GAS"
$ws.Cells.Item(2, 4).Value = 4.5
$ws.Cells.Item(2, 5).Value = 0.43
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1

# row 3: new placeholder NATURAL_GAS plant
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 99990300008
$ws.Cells.Item(3, 3).Value = "NATURAL_GAS"
$ws.Cells.Item(3, 4).Value = 4.2
$ws.Cells.Item(3, 5).Value = 0.61
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1

# row 4: previously HARD_COAL row (now moved up from old row 3)
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 19892800024
$ws.Cells.Item(4, 3).Value = "HARD_COAL"
$ws.Cells.Item(4, 4).Value = 3.5
$ws.Cells.Item(4, 5).Value = 0.33
$ws.Cells.Item(4, 6).Value = 24845.77
$ws.Cells.Item(4, 7).Value = 24845.77

# row 5: previously NUCLEAR row (now moved up from old row 6)
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 19851400028
$ws.Cells.Item(5, 3).Value = "NUCLEAR"
$ws.Cells.Item(5, 4).Value = 3.5
$ws.Cells.Item(5, 5).Value = 0.33
$ws.Cells.Item(5, 6).Value = 8599
$ws.Cells.Item(5, 7).Value = 8599

# ---------------------------------------------------------------------
# renewables: prepend three new placeholder plants, drop RunOfRiver,
# keep WindOn / OtherPV / WindOff shifted down.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("renewables")
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(7).Insert()

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 99992100002
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = "OtherPV"
$ws.Cells.Item(2, 6).Value = "-"
$ws.Cells.Item(2, 7).Value = "-"
$ws.Cells.Item(2, 8).Value = "-"
$ws.Cells.Item(2, 9).Value = "-"

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 99992400003
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 1.35
$ws.Cells.Item(3, 5).Value = "WindOn"
$ws.Cells.Item(3, 6).Value = "-"
$ws.Cells.Item(3, 7).Value = "-"
$ws.Cells.Item(3, 8).Value = "-"
$ws.Cells.Item(3, 9).Value = "-"

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 99992300007
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 2.7
$ws.Cells.Item(4, 5).Value = "WindOff"
$ws.Cells.Item(4, 6).Value = "-"
$ws.Cells.Item(4, 7).Value = "-"
$ws.Cells.Item(4, 8).Value = "-"
$ws.Cells.Item(4, 9).Value = "-"

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 20062400023
$ws.Cells.Item(5, 3).Value = 47547.50848700004
$ws.Cells.Item(5, 4).Value = 1.35
$ws.Cells.Item(5, 5).Value = "WindOn"
$ws.Cells.Item(5, 6).Value = "-"
$ws.Cells.Item(5, 7).Value = "-"
$ws.Cells.Item(5, 8).Value = "-"
$ws.Cells.Item(5, 9).Value = "-"

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 20102100030
$ws.Cells.Item(6, 3).Value = 53555.51607579708
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = "OtherPV"
$ws.Cells.Item(6, 6).Value = "-"
$ws.Cells.Item(6, 7).Value = "-"
$ws.Cells.Item(6, 8).Value = "-"
$ws.Cells.Item(6, 9).Value = "-"

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 20142300031
$ws.Cells.Item(7, 3).Value = 10271.8
$ws.Cells.Item(7, 4).Value = 2.7
$ws.Cells.Item(7, 5).Value = "WindOff"
$ws.Cells.Item(7, 6).Value = "-"
$ws.Cells.Item(7, 7).Value = "-"
$ws.Cells.Item(7, 8).Value = "-"
$ws.Cells.Item(7, 9).Value = "-"

# ---------------------------------------------------------------------
# storages: add the first (placeholder) storage plant row
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("storages")
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 99992600009
$ws.Cells.Item(2, 3).Value = "STORAGE"
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 1

# ---------------------------------------------------------------------
# biogas: prepend a new placeholder Biogas plant, shift the original down
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("biogas")
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 20000100021
$ws.Cells.Item(3, 3).Value = 4644.4034
$ws.Cells.Item(3, 4).Value = 1.9
$ws.Cells.Item(3, 5).Value = "Biogas"
$ws.Cells.Item(3, 6).Value = "-"
$ws.Cells.Item(3, 7).Value = "-"
$ws.Cells.Item(3, 8).Value = "-"
$ws.Cells.Item(3, 9).Value = "-"

$ws.Cells.Item(2, 2).Value = 99990100004
$ws.Cells.Item(2, 3).Value = 1

Write-Output "edit complete"
